# Generate Report for Handback
# Updates the zh-cn / de-de localization-status sheets to reflect a
# completed handback: status text, handback target/file/datetime columns,
# a hyperlink on the "Latest Target File" cell, and the wider columns
# that Excel uses once those columns hold real file-name content.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$sourceMd  = "31919597-49f0-4bd7-94a7-6977a4835286.md"
$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00c0ca81780699bacd0fbf4d09f6d52e5906f83a/e2e/31919597-49f0-4bd7-94a7-6977a4835286.md"
$statusText = "Handed back: in sync with en-US"

$zhXlf = "31919597-49f0-4bd7-94a7-6977a4835286.3dd97a080f5d38c4dcb8d2e3445955a81a426c68.zh-cn.xlf"
$deXlf = "31919597-49f0-4bd7-94a7-6977a4835286.3dd97a080f5d38c4dcb8d2e3445955a81a426c68.de-de.xlf"

# ---- zh-cn sheet ---------------------------------------------------
$ws2.Range("C2").Value2 = $statusText
$ws2.Range("J2").Value2 = $sourceMd
$f = $ws2.Range("J2").Font
$f.Underline = 2
$f.Color = 15570276
$ws2.Hyperlinks.Add($ws2.Range("J2"), $sourceUrl, [Type]::Missing, [Type]::Missing, $sourceMd)
$ws2.Range("K2").Value2 = $zhXlf
$ws2.Range("L2").Value2 = "2017-02-17 09:21:42"

$ws2.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$ws2.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664
$ws2.Range("K1").EntireColumn.ColumnWidth = 39.166666666666664

# ---- de-de sheet -----------------------------------------------------
$ws3.Range("C2").Value2 = $statusText
$ws3.Range("J2").Value2 = $sourceMd
$f3 = $ws3.Range("J2").Font
$f3.Underline = 2
$f3.Color = 15570276
$ws3.Hyperlinks.Add($ws3.Range("J2"), $sourceUrl, [Type]::Missing, [Type]::Missing, $sourceMd)
$ws3.Range("K2").Value2 = $deXlf
$ws3.Range("L2").Value2 = "2017-02-17 09:22:05"

$ws3.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$ws3.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664
$ws3.Range("K1").EntireColumn.ColumnWidth = 39.166666666666664

# ---- Overview sheet: zh-cn / de-de status columns widen in step ------
$ws1.Range("E1").EntireColumn.ColumnWidth = 29.166666666666668
$ws1.Range("F1").EntireColumn.ColumnWidth = 29.166666666666668
